$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Set data values for column I and J, rows 2-8
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 8

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 6

$ws.Range("I8").Value = 8
$ws.Range("J8").Value = 8
